$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Update the cached "datetimeFigureOut" date text (9/16/2018 -> 9/22/2018)
#    on every slide layout's Date Placeholder and on the slide master's Date
#    Placeholder. These placeholders are not present on individual slides in
#    this deck - they live on the master / layouts only.
# ---------------------------------------------------------------------------
$oldDate = "9/16/2018"
$newDate = "9/22/2018"

$master = $p.SlideMaster

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $txt = $shp.TextFrame.TextRange.Text
            if ($txt -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

for ($mi = 1; $mi -le $master.Shapes.Count; $mi++) {
    $mshp = $master.Shapes.Item($mi)
    if ($mshp.HasTextFrame -and $mshp.TextFrame.HasText) {
        $mtxt = $mshp.TextFrame.TextRange.Text
        if ($mtxt -eq $oldDate) {
            $mshp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Slide 11 ("Playing the game."): edit the sequence of bullet points.
#    - "The player rolls a dice and moves forward that many squares."
#      becomes "The player rolls " / "a die " / "and moves forward that many
#      squares." (the middle portion singled out as its own run).
#    - A new paragraph "Add one to the number of turns the player has
#      taken." is inserted right after that bullet.
#    - The old "Add one to the number of turns the player has taken."
#      bullet (which used to sit after "If they haven't won...") is removed.
# ---------------------------------------------------------------------------
$slide11 = $p.Slides.Item(11)
$body11 = $slide11.Shapes.Item(2).TextFrame.TextRange

[void]$body11.Replace("a dice ", "a die ")

$addLine = "Add one to the number of turns the player has taken."
$cr = [char]13

$dicePara = $body11.Paragraphs(2, 1)
[void]$dicePara.InsertAfter($cr + $addLine)

$paraTotal = ($body11.Text.ToCharArray() | Where-Object { $_ -eq $cr }).Count + 1
for ($i = 4; $i -le $paraTotal; $i++) {
    $candidate = $body11.Paragraphs($i, 1)
    if ($candidate.Text.TrimEnd($cr) -eq $addLine) {
        [void]$candidate.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 3. Slide 14 ("Starting to Get Statistics."): extend the sentence about
#    restructuring example5.
# ---------------------------------------------------------------------------
$slide14 = $p.Slides.Item(14)
$body14 = $slide14.Shapes.Item(2).TextFrame.TextRange
[void]$body14.Replace(
    "Restructure example5 so it plays the game ten thousand ",
    "Restructure example5 so it plays the game ten thousand times "
)
